$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 15:55"

# 2. Swap country-name labels where the ranking order changed
#    (row stays the same, but which country occupies that row swaps)
$ws.Range("A61").Value = "Serbia"
$ws.Range("A62").Value = "Irlanda"
$ws.Range("A79").Value = "Estado de Palestina"
$ws.Range("A80").Value = "Bosnia y Herzegovina"
$ws.Range("A96").Value = "Zambia"
$ws.Range("A97").Value = "Mauritania"

# 3. Updated case-count figures
$ws.Range("B4").Value = 4766323
$ws.Range("C4").Value = 2005
$ws.Range("E4").Value = 2245234
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 157924
$ws.Range("B6").Value = 1767836
$ws.Range("C6").Value = 15917
$ws.Range("D6").Value = 1157139
$ws.Range("E6").Value = 573127
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 37570
$ws.Range("B17").Value = 278835
$ws.Range("C17").Value = 1357
$ws.Range("D17").Value = 240081
$ws.Range("E17").Value = 35837
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 2917
$ws.Range("B21").Value = 211208
$ws.Range("C21").Value = 131
$ws.Range("E21").Value = 8382
$ws.Range("B44").Value = 55098
$ws.Range("C44").Value = 366
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 6149
$ws.Range("B61").Value = 26193
$ws.Range("C61").Value = 311
$ws.Range("D61").Value = 14047
$ws.Range("E61").Value = 11556
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 590
$ws.Range("B62").Value = 26109
$ws.Range("D62").Value = 23364
$ws.Range("E62").Value = 982
$ws.Range("H62").Value = 1763
$ws.Range("B66").Value = 22053
$ws.Range("C66").Value = 690
$ws.Range("E66").Value = 13265
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 369
$ws.Range("B72").Value = 17923
$ws.Range("C72").Value = 641
$ws.Range("D72").Value = 10420
$ws.Range("E72").Value = 7295
$ws.Range("B79").Value = 12297
$ws.Range("C79").Value = 137
$ws.Range("D79").Value = 5390
$ws.Range("E79").Value = 6823
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 84
$ws.Range("B80").Value = 12296
$ws.Range("C80").Value = 420
$ws.Range("D80").Value = 6312
$ws.Range("E80").Value = 5632
$ws.Range("G80").Value = 13
$ws.Range("H80").Value = 352
$ws.Range("B86").Value = 9263
$ws.Range("C86").Value = 10
$ws.Range("E86").Value = 256
$ws.Range("B94").Value = 7317
$ws.Range("C94").Value = 9
$ws.Range("D94").Value = 6480
$ws.Range("E94").Value = 791
$ws.Range("B96").Value = 6347
$ws.Range("C96").Value = 119
$ws.Range("D96").Value = 4493
$ws.Range("E96").Value = 1684
$ws.Range("G96").Value = 5
$ws.Range("H96").Value = 170
$ws.Range("B97").Value = 6319
$ws.Range("D97").Value = 5043
$ws.Range("E97").Value = 1119
$ws.Range("H97").Value = 157
$ws.Range("B101").Value = 5161
$ws.Range("C101").Value = 77
$ws.Range("D101").Value = 5019
$ws.Range("E101").Value = 83
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 59
$ws.Range("B118").Value = 2817
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 292
$ws.Range("D146").Value = 945
$ws.Range("E146").Value = 145
$ws.Range("D147").Value = 856
$ws.Range("E147").Value = 264
$ws.Range("E161").Value = 241
$ws.Range("G161").Value = 3
$ws.Range("H161").Value = 6
$ws.Range("D169").Value = 298
$ws.Range("E169").Value = 49
$ws.Range("D179").Value = 183
$ws.Range("E179").Value = 5
